$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games), row-by-row B:E and G (F unchanged)
$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    3 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 29.84159230404497)
    4 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286)
    5 = @(0.3048080303191223, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.626907116734944)
    6 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    7 = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 31.61296591696135)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G
}
